$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 268, shifting existing rows 268:354 down to 269:355
$ws.Rows.Item(268).Insert()

# Populate the newly inserted row 268 with the new data record
$ws.Cells.Item(268, 1).Value = 7
$ws.Cells.Item(268, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(268, 3).Value = "Ñuble"
$ws.Cells.Item(268, 4).Value = 44988
$ws.Cells.Item(268, 5).Value = 16
$ws.Cells.Item(268, 6).Value = 100112003
$ws.Cells.Item(268, 7).Value = "Ajo"
$ws.Cells.Item(268, 8).Value = "Chino"
$ws.Cells.Item(268, 9).Value = "Primera"
$ws.Cells.Item(268, 10).Value = 50
$ws.Cells.Item(268, 11).Value = 21000
$ws.Cells.Item(268, 12).Value = 21000
$ws.Cells.Item(268, 13).Value = 21000
$ws.Cells.Item(268, 14).Value = "$/malla 10 kilos"
$ws.Cells.Item(268, 15).Value = "China"
$ws.Cells.Item(268, 16).Value = 2100
$ws.Cells.Item(268, 17).Value = 10
$ws.Cells.Item(268, 18).Value = "Hortaliza"
